$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.231.95'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').Value = '1.820.87'
$ws.Range('E3').Value = '  -1.90%  '
$ws.Range('E4').Value = '  -1.30%  '
$ws.Range('D5').Value = '''314.28'
$ws.Range('E5').Value = '  -1.79%  '
$ws.Range('D6').Value = '''1.004'
$ws.Range('E7').Value = '  -2.19%  '
$ws.Range('D8').Value = '''0.3669'
$ws.Range('E8').Value = '  -2.83%  '
$ws.Range('D9').Value = '''45.86'
$ws.Range('E9').Value = '  -1.72%  '
$ws.Range('D10').Value = '''0.07209'
$ws.Range('E10').Value = '  -2.60%  '
$ws.Range('D11').Value = '''0.8594'
$ws.Range('E11').Value = '  -2.47%  '
$ws.Range('D12').Value = '''20.95'
$ws.Range('E12').Value = '  -2.94%  '
$ws.Range('D13').Value = '1.832.41'
$ws.Range('E13').Value = '  -1.45%  '
$ws.Range('D14').Value = '''6.648'
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('D15').Value = '''0.07109'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').Value = '''5.298'
$ws.Range('E16').Value = '  -3.35%  '
$ws.Range('D17').Value = '''87.79'
$ws.Range('E17').Value = '  +1.01%  '
$ws.Range('E18').Value = '  -1.53%  '
$ws.Range('D19').Value = '''0.000008834'
$ws.Range('E19').Value = '  -2.36%  '
$ws.Range('D20').Value = '''1.004'
$ws.Range('E20').Value = '  -1.23%  '
$ws.Range('D21').Value = '''15.02'
$ws.Range('E21').Value = '  -2.32%  '
$ws.Range('D22').Value = '27.257.73'
$ws.Range('E22').Value = '  -1.57%  '
$ws.Range('D23').Value = '''5.129'
$ws.Range('E23').Value = '  -2.84%  '
$ws.Range('D24').Value = '''10.87'
$ws.Range('E24').Value = '  -2.26%  '
$ws.Range('D25').Value = '2.055.99'
$ws.Range('E25').Value = '  -1.50%  '
$ws.Range('D26').Value = '''1.999'
$ws.Range('E26').Value = '  -1.64%  '
$ws.Range('D27').Value = '''153.14'
$ws.Range('E27').Value = '  -2.28%  '
$ws.Range('D28').Value = '''18.28'
$ws.Range('E28').Value = '  -2.02%  '
$ws.Range('E29').Value = '  +5.61%  '
$ws.Range('D30').Value = '''5.220'
$ws.Range('E30').Value = '  -2.77%  '
$ws.Range('D31').Value = '''116.07'
$ws.Range('E31').Value = '  -3.96%  '
$ws.Range('D32').Value = '''0.08884'
$ws.Range('E32').Value = '  -1.77%  '
$ws.Range('D33').Value = '''0.7594'
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('D34').Value = '''1.190'
$ws.Range('E34').Value = '  -2.14%  '
$ws.Range('D35').Value = '''4.452'
$ws.Range('E35').Value = '  -2.12%  '
$ws.Range('D36').Value = '''2.817'
$ws.Range('E36').Value = '  -6.95%  '
$ws.Range('D37').Value = '''1.005'
$ws.Range('E37').Value = '  -1.28%  '
$ws.Range('D38').Value = '''1.111'
$ws.Range('E38').Value = '  -2.44%  '
$ws.Range('D39').Value = '''0.01954'
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('D40').Value = '''0.05247'
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('D41').Value = '''2.896'
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('D42').Value = '''7.041'
$ws.Range('E42').Value = '  +1.30%  '
$ws.Range('D43').Value = '''0.1675'
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').Value = '''0.5011'
$ws.Range('E44').Value = '  -3.34%  '
$ws.Range('D45').Value = '''8.604'
$ws.Range('E45').Value = '  -0.97%  '
$ws.Range('E46').Value = '  -2.06%  '
$ws.Range('D47').Value = '''106.56'
$ws.Range('E47').Value = '  -3.13%  '
$ws.Range('D48').Value = '''0.4686'
$ws.Range('E48').Value = '  -0.67%  '
$ws.Range('D49').Value = '''1.004'
$ws.Range('E49').Value = '  -1.32%  '
$ws.Range('D50').Value = '''0.06391'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('D51').Value = '''1.657'
$ws.Range('E51').Value = '  -3.03%  '
